$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.703862428665161
$ws.Range("B1").Value = 1.854307293891907
$ws.Range("C1").Value = 1.847112655639648
$ws.Range("D1").Value = 2.405254364013672
$ws.Range("E1").Value = 3.022364616394043
